# before.xlsx col F ("柴油产销率") and col G ("柴油销售量") are dropped, and within each
# year the "B"-labelled and "C"-labelled quarter rows are swapped (A/D/E and the
# applicable B/C figures move with them) so rows read A, C, B, D per year.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# swap row 3 <-> row 4
$ws.Range('A3').Value2 = '2000年C'
$ws.Range('B3').Value2 = 99.3
$ws.Range('D3').Value2 = 60.9
$ws.Range('E3').Value2 = 5226.9
$ws.Range('A4').Value2 = '2000年B'
$ws.Range('B4').Value2 = 99.5
$ws.Range('D4').Value2 = 20.8
$ws.Range('E4').Value2 = 3288.6

# swap row 7 <-> row 8
$ws.Range('A7').Value2 = '2001年C'
$ws.Range('B7').Value2 = 99
$ws.Range('C7').Value2 = -0.3
$ws.Range('D7').Value2 = 28
$ws.Range('E7').Value2 = 5456.1
$ws.Range('A8').Value2 = '2001年B'
$ws.Range('B8').Value2 = 100
$ws.Range('C8').ClearContents()
$ws.Range('D8').Value2 = 15.4
$ws.Range('E8').Value2 = 3636.3

# swap row 11 <-> row 12
$ws.Range('A11').Value2 = '2002年C'
$ws.Range('B11').Value2 = 100.8
$ws.Range('D11').Value2 = -7.1
$ws.Range('E11').Value2 = 5705.3
$ws.Range('A12').Value2 = '2002年B'
$ws.Range('B12').Value2 = 99.59999999999999
$ws.Range('D12').Value2 = 7.4
$ws.Range('E12').Value2 = 3785.5

# swap row 15 <-> row 16
$ws.Range('A15').Value2 = '2003年C'
$ws.Range('B15').Value2 = 100.4
$ws.Range('C15').Value2 = -0.4
$ws.Range('D15').Value2 = -9.199999999999999
$ws.Range('E15').Value2 = 6247.4
$ws.Range('A16').Value2 = '2003年B'
$ws.Range('B16').Value2 = 99.59999999999999
$ws.Range('C16').Value2 = 0
$ws.Range('D16').Value2 = -1
$ws.Range('E16').Value2 = 4026.3

# swap row 19 <-> row 20
$ws.Range('A19').Value2 = '2004年C'
$ws.Range('B19').Value2 = 99.7
$ws.Range('C19').Value2 = -0.7
$ws.Range('D19').Value2 = 24.7
$ws.Range('E19').Value2 = 7472.7
$ws.Range('A20').Value2 = '2004年B'
$ws.Range('B20').Value2 = 99.8
$ws.Range('C20').Value2 = 0.2
$ws.Range('D20').Value2 = 12.6
$ws.Range('E20').Value2 = 4899

# swap row 23 <-> row 24
$ws.Range('A23').Value2 = '2005年C'
$ws.Range('B23').Value2 = 99.8
$ws.Range('C23').Value2 = 0.1
$ws.Range('D23').Value2 = 1.9
$ws.Range('E23').Value2 = 8130.7
$ws.Range('A24').Value2 = '2005年B'
$ws.Range('B24').Value2 = 100.3
$ws.Range('C24').Value2 = 0.5
$ws.Range('D24').Value2 = -12.1
$ws.Range('E24').Value2 = 5348.5

# swap row 27 <-> row 28
$ws.Range('A27').Value2 = '2006年C'
$ws.Range('B27').Value2 = 99.40000000000001
$ws.Range('C27').Value2 = -0.4
$ws.Range('D27').Value2 = 40.9
$ws.Range('E27').Value2 = 8574.299999999999
$ws.Range('A28').Value2 = '2006年B'
$ws.Range('B28').Value2 = 99.5
$ws.Range('C28').Value2 = -0.8
$ws.Range('D28').Value2 = 21.4
$ws.Range('E28').Value2 = 5675.3

# swap row 31 <-> row 32
$ws.Range('A31').Value2 = '2007年C'
$ws.Range('C31').Value2 = 0.4
$ws.Range('D31').Value2 = 10.2
$ws.Range('E31').Value2 = 9128.5
$ws.Range('A32').Value2 = '2007年B'
$ws.Range('C32').Value2 = 0.3
$ws.Range('D32').Value2 = 7.6
$ws.Range('E32').Value2 = 5932.9

# swap row 35 <-> row 36
$ws.Range('A35').Value2 = '2008年C'
$ws.Range('B35').Value2 = 99.5
$ws.Range('D35').Value2 = 32
$ws.Range('E35').Value2 = 9984
$ws.Range('A36').Value2 = '2008年B'
$ws.Range('B36').Value2 = 99.59999999999999
$ws.Range('D36').Value2 = 17.6
$ws.Range('E36').Value2 = 6502.5

# swap row 39 <-> row 40
$ws.Range('A39').Value2 = '2009年C'
$ws.Range('B39').Value2 = 99.3
$ws.Range('C39').Value2 = -0.3
$ws.Range('D39').Value2 = 40.7
$ws.Range('E39').Value2 = 10064.9
$ws.Range('A40').Value2 = '2009年B'
$ws.Range('B40').Value2 = 99.59999999999999
$ws.Range('C40').Value2 = 0.1
$ws.Range('D40').Value2 = 16.2
$ws.Range('E40').Value2 = 6331.8

# swap row 43 <-> row 44
$ws.Range('A43').Value2 = '2010年C'
$ws.Range('B43').Value2 = 100.3
$ws.Range('C43').Value2 = 1.2
$ws.Range('D43').Value2 = -14.4
$ws.Range('E43').Value2 = 11528.1
$ws.Range('A44').Value2 = '2010年B'
$ws.Range('B44').Value2 = 100.4
$ws.Range('C44').Value2 = 0.9
$ws.Range('D44').Value2 = -11.8
$ws.Range('E44').Value2 = 7570

# swap row 47 <-> row 48
$ws.Range('A47').Value2 = '2011年C'
$ws.Range('B47').Value2 = 99.8
$ws.Range('C47').Value2 = -0.5
$ws.Range('D47').Value2 = 13.5
$ws.Range('E47').Value2 = 12397.4
$ws.Range('A48').Value2 = '2011年B'
$ws.Range('B48').Value2 = 100
$ws.Range('C48').Value2 = -0.1
$ws.Range('D48').Value2 = 0
$ws.Range('E48').Value2 = 8309.4

# swap row 51 <-> row 52
$ws.Range('A51').Value2 = '2012年C'
$ws.Range('B51').Value2 = 100
$ws.Range('C51').Value2 = 0.1
$ws.Range('D51').Value2 = 1.7
$ws.Range('E51').Value2 = 12619.4
$ws.Range('A52').Value2 = '2012年B'
$ws.Range('B52').Value2 = 99.2
$ws.Range('C52').Value2 = -0.8
$ws.Range('D52').Value2 = 27.1
$ws.Range('E52').Value2 = 8411

# swap row 55 <-> row 56
$ws.Range('A55').Value2 = '2013年C'
$ws.Range('B55').Value2 = 99.40000000000001
$ws.Range('C55').Value2 = -0.4
$ws.Range('D55').Value2 = 28
$ws.Range('E55').Value2 = 12762.6
$ws.Range('A56').Value2 = '2013年B'
$ws.Range('B56').Value2 = 98.8
$ws.Range('C56').Value2 = -0.2
$ws.Range('D56').Value2 = 39
$ws.Range('E56').Value2 = 8454.9

# swap row 59 <-> row 60
$ws.Range('A59').Value2 = '2014年C'
$ws.Range('C59').Value2 = -0.4
$ws.Range('D59').Value2 = 28.6
$ws.Range('E59').Value2 = 12810.2
$ws.Range('A60').Value2 = '2014年B'
$ws.Range('C60').Value2 = 0.1
$ws.Range('D60').Value2 = 20.2
$ws.Range('E60').Value2 = 8452.799999999999

# swap row 63 <-> row 64
$ws.Range('A63').Value2 = '2015年C'
$ws.Range('B63').Value2 = 99.09999999999999
$ws.Range('C63').Value2 = 0.1
$ws.Range('D63').Value2 = 18.9
$ws.Range('E63').Value2 = 13212.6
$ws.Range('A64').Value2 = '2015年B'
$ws.Range('B64').Value2 = 98.59999999999999
$ws.Range('C64').Value2 = -0.3
$ws.Range('D64').Value2 = 29
$ws.Range('E64').Value2 = 8767.9

# 柴油产销率 (F) and 柴油销售量 (G) columns are removed entirely (not just cleared),
# which also shrinks the sheet dimension from A1:G64 to A1:E64.
$ws.Range('F1:G64').Delete()
